$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

# Update column B (Taxonsorteringsordning) from 57069 to 57073
# for every data row where it currently holds the old value.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    if ($cell.Value2 -eq 57069) {
        $cell.Value = 57073
    }
}
